$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '98.962.54'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  +1.95%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.371.45'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  +7.81%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.04%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '259.13'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  +8.84%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '629.49'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  +3.01%  '
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  +25.53%  '
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.393'
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  +2.28%  '
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -0.09%  '
$r.Style = "Normal"
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '0.882'
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  +12.90%  '
$r.Style = "Normal"
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '3.368.43'
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  +7.83%  '
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  +0.90%  '
$r.Style = "Normal"
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '98.772.13'
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  +1.97%  '
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '36.13'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  +7.23%  '
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  +3.28%  '
$r.Style = "Normal"
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '3.961.84'
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  +7.17%  '
$r.Style = "Normal"
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '5.53'
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  +1.85%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '3.369.05'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  +7.94%  '
$r.Style = "Normal"
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '3.58'
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  +1.52%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '15.28'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  +5.26%  '
$r.Style = "Normal"
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '494.94'
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -5.63%  '
$r.Style = "Normal"
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '6.19'
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  +8.88%  '
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  +9.92%  '
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '9.48'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  +7.40%  '
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '5.71'
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  +4.55%  '
$r.Style = "Normal"
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '89.91'
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  +1.78%  '
$r.Style = "Normal"
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '12.00'
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  +3.81%  '
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  +7.35%  '
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  +19.98%  '
$r.Style = "Normal"
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  +0.05%  '
$r.Style = "Normal"
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '0.191'
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  +9.92%  '
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  +8.57%  '
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '0.998'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  +15.19%  '
$r.Style = "Normal"
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '9.59'
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  +6.49%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '27.86'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  +3.05%  '
$r.Style = "Normal"
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '0.152'
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +0.28%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '7.36'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  +1.37%  '
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  +5.62%  '
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  +5.60%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '498.23'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  +4.30%  '
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  +2.72%  '
$r.Style = "Normal"
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '3.88'
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  +7.47%  '
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  +3.58%  '
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '3.30'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  +4.73%  '
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '0.786'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  +12.62%  '
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  +0.02%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '160.61'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  -0.89%  '
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  +2.39%  '
$r.Style = "Normal"
$r = $ws.Range('B49')
$r.NumberFormat = "@"
$r.Value = 'Mantle'
$r.Style = "Normal"
$r = $ws.Range('C49')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$r.Style = "Normal"
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '0.832'
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  +14.07%  '
$r.Style = "Normal"
$r = $ws.Range('B50')
$r.NumberFormat = "@"
$r.Value = 'Filecoin'
$r.Style = "Normal"
$r = $ws.Range('C50')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '4.67'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  +3.96%  '
$r.Style = "Normal"
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '46.13'
$r.Style = "Normal"
